$wb = $excel.ActiveWorkbook

# --- Sheet "Moorings": fix the recorded longitude for GA03FLMA-00001 ---
$moorings = $wb.Worksheets.Item("Moorings")
$moorings.Range("I2").Value = "42° 52.24' W"
$moorings.Range("L2").Value = "Wrong anchor position recovered in Argentine-1 Cruise report. Fixed here"

# --- Sheet "Asset_Cal_Info": update every CC_lon calibration value that used the old longitude ---
$assetCal = $wb.Worksheets.Item("Asset_Cal_Info")
$lonRows = @(22,26,28,37,42,47,52,57,62,67,72,77,82,87,92)
foreach ($r in $lonRows) {
    $assetCal.Range("H$r").Value = -42.870666666666665
}

$wb.Application.Calculate()

# --- cosmetic: restore the view/selection state recorded in the workbook ---
$moorings.Activate()
$moorings.Range("N2").Select()
$assetCal.Activate()
$assetCal.Application.ActiveWindow.ScrollRow = 1
$assetCal.Range("O32").Select()
